# Weekly update to the "Zapallo italiano" price sheet:
# a new daily record is inserted at row 335 (pushing the existing
# rows 335-385 down to 336-386) and populated with the latest data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 335; everything below shifts down by one.
$ws.Rows.Item(335).Insert()

# Fill in the new row 335 with the new record's data.
$ws.Cells.Item(335, 1).Value2  = 4
$ws.Cells.Item(335, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(335, 3).Value2  = "Los Lagos"
$ws.Cells.Item(335, 4).Value2  = 45034
$ws.Cells.Item(335, 5).Value2  = 10
$ws.Cells.Item(335, 6).Value2  = 100112032
$ws.Cells.Item(335, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(335, 8).Value2  = "Sin especificar"
$ws.Cells.Item(335, 9).Value2  = "Primera"
$ws.Cells.Item(335, 10).Value2 = 250
$ws.Cells.Item(335, 11).Value2 = 13000
$ws.Cells.Item(335, 12).Value2 = 13000
$ws.Cells.Item(335, 13).Value2 = 13000
$ws.Cells.Item(335, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(335, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(335, 16).Value2 = 260
$ws.Cells.Item(335, 17).Value2 = 50
$ws.Cells.Item(335, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by
# the rest of column D.
$ws.Cells.Item(335, 4).NumberFormat = $ws.Cells.Item(336, 4).NumberFormat
